# Update for final supersite list:
#  - Erie Elementary -> Meadowlark (rows 42-50)
#  - Southern Hills MS -> Fairview HS (various rows)
#  - A handful of single-row supersite values moved out of the middle of the
#    sheet and appended as new rows at the bottom (195-205)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Erie Elementary -> Meadowlark -----------------------------------
$erieRows = 42,43,44,45,46,47,48,49,50
foreach ($r in $erieRows) {
    $ws.Range("H$r").Value = "Meadowlark"
}

# --- Southern Hills MS -> Fairview HS ---------------------------------
$fairviewRows = 150,163,164,165,166,167,170,171,172,174,176,177,178,184
foreach ($r in $fairviewRows) {
    $ws.Range("H$r").Value = "Fairview HS"
}

# --- Clear single supersite entries that move down to the new rows ----
$clearRows = 21,34,52,78,94,98,121,135,146,173,175
foreach ($r in $clearRows) {
    $ws.Range("H$r").ClearContents()
}

# --- Append the new rows 195-205 with only column H populated ---------
$newRowValues = @{
    195 = "Casey MS"
    196 = "Centennial MS"
    197 = "Fairview HS"
    198 = "Fairview HS"
    199 = "Niwot HS"
    200 = "Centaurus HS"
    201 = "Centaurus HS"
    202 = "Altona MS"
    203 = "Longs Peak MS"
    204 = "Timberline K-8"
    205 = "Louisville MS"
}

foreach ($r in 195..205) {
    $ws.Range("H$r").Value = $newRowValues[$r]
}
